$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98

# Column A holds a date-shaped string ("2025-10-24"). Writing it plainly would
# cause Excel to auto-convert it to a date serial number, which is not what the
# source data (and the rest of the column) represents - every existing cell in
# this sheet stores plain text. Prefixing with a quote forces text entry, then
# ClearFormats() removes the leftover "quote prefix" cell format so the cell
# ends up with the same (default) style as all the other rows.
$ws.Cells.Item($row, 1).Value = "'2025-10-24"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "21:21:27"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,794.4801"
